# Update cryptocurrency price/volume snapshot on Sheet1 (rows 2-51).
# For Price (column D) values that look like plain numbers, force a text
# number format first so Excel stores them as text (matching the source
# data, e.g. thousand-dot formatted prices) instead of auto-converting
# them to numeric values; then reset the cell style back to Normal so no
# stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.279.43'
$ws.Cells.Item(2, 5).Value = '  -0.16%  '

$ws.Cells.Item(3, 4).Value = '2.644.86'
$ws.Cells.Item(3, 5).Value = '  +0.41%  '

$ws.Cells.Item(4, 5).Value = '  +0.12%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '597.74'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.13%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '155.45'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.81%  '

$ws.Cells.Item(7, 5).Value = '  +0.03%  '

$ws.Cells.Item(8, 5).Value = '  -0.05%  '

$ws.Cells.Item(9, 5).Value = '  +7.36%  '

$ws.Cells.Item(10, 5).Value = '  -0.93%  '

$ws.Cells.Item(11, 5).Value = '  +0.40%  '

$ws.Cells.Item(12, 5).Value = '  +1.37%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '28.05'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +1.63%  '

$ws.Cells.Item(14, 5).Value = '  +2.33%  '

$ws.Cells.Item(15, 4).Value = '3.127.62'
$ws.Cells.Item(15, 5).Value = '  +0.59%  '

$ws.Cells.Item(16, 4).Value = '68.223.32'
$ws.Cells.Item(16, 5).Value = '  -0.02%  '

$ws.Cells.Item(17, 4).Value = '2.637.92'
$ws.Cells.Item(17, 5).Value = '  -0.80%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '11.40'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.06%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '363.96'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -1.46%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.46'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.84%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.39'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +3.43%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.83'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.20%  '

$ws.Cells.Item(23, 5).Value = '  -0.89%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '75.32'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +3.32%  '

$ws.Cells.Item(25, 5).Value = '  -0.09%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '9.68'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -2.08%  '

$ws.Cells.Item(27, 5).Value = '  +1.85%  '

$ws.Cells.Item(28, 4).Value = '2.782.29'
$ws.Cells.Item(28, 5).Value = '  +0.61%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.998'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.35%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '560.22'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -2.19%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '8.03'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +0.99%  '

$ws.Cells.Item(32, 5).Value = '  +0.98%  '

$ws.Cells.Item(33, 5).Value = '  +0.62%  '

$ws.Cells.Item(34, 5).Value = '  +1.72%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.14%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.56'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +3.28%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '161.17'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.56%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '19.36'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +1.09%  '

$ws.Cells.Item(39, 5).Value = '  +1.64%  '

$ws.Cells.Item(40, 5).Value = '  -0.87%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '5.33'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -0.23%  '

$ws.Cells.Item(42, 4).Value = '0.0₆0341'
$ws.Cells.Item(42, 5).Value = '  +4.53%  '

$ws.Cells.Item(43, 5).Value = '  -0.30%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '17.78'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.03%  '

$ws.Cells.Item(45, 5).Value = '  +0.01%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '40.44'
$ws.Cells.Item(46, 4).Style = "Normal"

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '159.22'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +2.67%  '

$ws.Cells.Item(48, 5).Value = '  +1.26%  '

$ws.Cells.Item(49, 5).Value = '  +0.32%  '

$ws.Cells.Item(50, 2).Value = 'Optimism'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.69'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -0.09%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0785'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.97%  '
